# Update countries & provincias Spain
# Refresh the COVID-19 snapshot in the "Pais" sheet: new totals per
# country (re-sorted descending by "Casos totales", which reshuffles a
# few rows), plus the "updated as of" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot timestamp
$ws.Range("A1").Value = 'Datos actualizados a 2 de Abril de 2020 a las 19:20'

# Row 4
$ws.Range("B4").Value = 235972
$ws.Range("C4").Value = 20969
$ws.Range("D4").Value = 10304
$ws.Range("E4").Value = 219893
$ws.Range("G4").Value = 673
$ws.Range("H4").Value = 5775

# Row 6
$ws.Range("E6").Value = 73399
$ws.Range("G6").Value = 709
$ws.Range("H6").Value = 10096

# Row 7
$ws.Range("B7").Value = 83459
$ws.Range("C7").Value = 5478
$ws.Range("D7").Value = 21400
$ws.Range("E7").Value = 61011
$ws.Range("G7").Value = 117
$ws.Range("H7").Value = 1048

# Row 16
$ws.Range("B16").Value = 11076
$ws.Range("C16").Value = 365
$ws.Range("E16").Value = 9169

# Row 17
$ws.Range("B17").Value = 11068
$ws.Range("C17").Value = 1337
$ws.Range("E17").Value = 9043
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 134

# Row 20
$ws.Range("B20").Value = 7031
$ws.Range("C20").Value = 151
$ws.Range("E20").Value = 6652

# Row 24
$ws.Range("B24").Value = 5128
$ws.Range("C24").Value = 251
$ws.Range("E24").Value = 5046
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 50

# Row 25
$ws.Range("A25").Value = 'Irlanda'
$ws.Range("B25").Value = 3849
$ws.Range("C25").Value = 402
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 3746
$ws.Range("F25").Value = 109
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 98

# Row 26
$ws.Range("A26").Value = 'Chequia'
$ws.Range("B26").Value = 3805
$ws.Range("C26").Value = 216
$ws.Range("D26").Value = 67
$ws.Range("E26").Value = 3694
$ws.Range("F26").Value = 72
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 44

# Row 27
$ws.Range("A27").Value = 'Rusia'
$ws.Range("B27").Value = 3548
$ws.Range("C27").Value = 771
$ws.Range("D27").Value = 235
$ws.Range("E27").Value = 3283
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 30

# Row 29
$ws.Range("B29").Value = 3386
$ws.Range("C29").Value = 279
$ws.Range("E29").Value = 2174

# Row 36
$ws.Range("E36").Value = 2246
$ws.Range("F36").Value = 9
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 33

# Row 58
$ws.Range("A58").Value = 'Egipto'
$ws.Range("B58").Value = 850
$ws.Range("C58").Value = 71
$ws.Range("D58").Value = 179
$ws.Range("E58").Value = 619
$ws.Range("F58").Value = 0
$ws.Range("H58").Value = 52

# Row 59
$ws.Range("A59").Value = 'Catar'
$ws.Range("B59").Value = 835
$ws.Range("D59").Value = 71
$ws.Range("E59").Value = 762
$ws.Range("F59").Value = 37
$ws.Range("H59").Value = 2

# Row 60
$ws.Range("A60").Value = 'Emiratos Arabes Unidos'
$ws.Range("B60").Value = 814
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 61
$ws.Range("E60").Value = 745
$ws.Range("F60").Value = 2
$ws.Range("H60").Value = 8

# Row 61
$ws.Range("A61").Value = 'Ucrania'
$ws.Range("B61").Value = 804
$ws.Range("C61").Value = 10
$ws.Range("D61").Value = 13
$ws.Range("E61").Value = 771
$ws.Range("F61").Value = 0
$ws.Range("H61").Value = 20

# Row 62
$ws.Range("A62").Value = 'Hong Kong'
$ws.Range("B62").Value = 802
$ws.Range("C62").Value = 36
$ws.Range("D62").Value = 154
$ws.Range("E62").Value = 644
$ws.Range("F62").Value = 8
$ws.Range("H62").Value = 4

# Row 63
$ws.Range("A63").Value = 'Nueva Zelanda'
$ws.Range("B63").Value = 797
$ws.Range("C63").Value = 89
$ws.Range("D63").Value = 92
$ws.Range("E63").Value = 704
$ws.Range("F63").Value = 2
$ws.Range("H63").Value = 1

# Row 66
$ws.Range("B66").Value = 691
$ws.Range("C66").Value = 37
$ws.Range("D66").Value = 30
$ws.Range("E66").Value = 617
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 44

# Row 89
$ws.Range("A89").Value = 'Jordania'
$ws.Range("B89").Value = 299
$ws.Range("C89").Value = 21
$ws.Range("D89").Value = 45
$ws.Range("E89").Value = 249
$ws.Range("F89").Value = 5
$ws.Range("H89").Value = 5

# Row 90
$ws.Range("A90").Value = 'Burkina Faso'
$ws.Range("B90").Value = 288
$ws.Range("C90").Value = 6
$ws.Range("D90").Value = 50
$ws.Range("E90").Value = 222
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 16

# Row 91
$ws.Range("A91").Value = 'Camerun'
$ws.Range("B91").Value = 284
$ws.Range("C91").Value = 51
$ws.Range("D91").Value = 10
$ws.Range("E91").Value = 267
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 7

# Row 93
$ws.Range("B93").Value = 258
$ws.Range("C93").Value = 21
$ws.Range("E93").Value = 244

# Row 149
$ws.Range("A149").Value = 'Eritrea'
$ws.Range("C149").Value = 7
$ws.Range("E149").Value = 22
$ws.Range("H149").Value = 0

# Row 150
$ws.Range("A150").Value = 'Islas Caimanes'
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 21

# Row 151
$ws.Range("A151").Value = 'San Martin (Parte Francesa)'
$ws.Range("E151").Value = 19
$ws.Range("H151").Value = 1

# Row 152
$ws.Range("A152").Value = 'Congo'
$ws.Range("B152").Value = 22
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 18
$ws.Range("H152").Value = 2

# Row 153
$ws.Range("A153").Value = 'Bahamas'
$ws.Range("B153").Value = 21
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 1

# Row 154
$ws.Range("A154").Value = 'Birmania'
$ws.Range("C154").Value = 4
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 19

# Row 155
$ws.Range("A155").Value = 'Tanzania'
$ws.Range("B155").Value = 20
$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 17
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 1

# Row 156
$ws.Range("A156").Value = 'Guyana'
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 15
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 4

# Row 157
$ws.Range("A157").Value = 'Maldivas'
$ws.Range("B157").Value = 19
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 13
$ws.Range("E157").Value = 6

